$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the issue id "10", which looks numeric. Force it to be
# stored as text (matching the rest of the sheet, where every value -
# including other numeric-looking ids - is text) by setting a text
# number format before assigning the value.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "10"

$ws.Range("B6").Value = "[BUG] ghdwedgweyd"
$ws.Range("C6").Value = "open"
$ws.Range("D6").Value = "2025-03-24T08:32:16Z"
$ws.Range("E6").Value = "bug"
